$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 46 (shifts existing rows 46-165 down to 47-166)
$ws.Rows("46").Insert()

# Populate the newly inserted row 46 with a new weekly record.
# All fixed/contextual columns mirror the record that is now in row 47
# (the former row 46), except for the Fecha (D) and Volumen (J) values,
# which are the new data point's own values.
$ws.Range("A46").Value = 4
$ws.Range("B46").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C46").Value = "Los Lagos"
$ws.Range("D46").Value = 44526
$ws.Range("E46").Value = 10
$ws.Range("F46").Value = 100112017
$ws.Range("G46").Value = "Apio"
$ws.Range("H46").Value = "Americana (o)"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 35
$ws.Range("K46").Value = 11000
$ws.Range("L46").Value = 11000
$ws.Range("M46").Value = 11000
$ws.Range("N46").Value = "$/docena de matas"
$ws.Range("O46").Value = "Región de Coquimbo"
$ws.Range("P46").Value = 1833
$ws.Range("Q46").Value = 6
$ws.Range("R46").Value = "Hortaliza"

# Match the date style used by the rest of the Fecha column
$ws.Range("D46").NumberFormat = $ws.Range("D47").NumberFormat
